$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell, preserving text representation
# for values that look numeric (so they stay shared strings, not numbers),
# matching the source data which stores these as text.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: 이노그리드
$ws.Range("A2").Value = "이노그리드"
$ws.Range("B2").Value = "2024.05.31~06.07"
$ws.Range("C2").Value = "29,000~35,000"
$ws.Range("D2").Value = "-"
Set-TextValue $ws.Range("E2") "17400"
$ws.Range("F2").Value = "한국투자증권"

# Row 3: 씨어스테크놀로지
$ws.Range("A3").Value = "씨어스테크놀로지"
$ws.Range("B3").Value = "2024.05.23~05.29"
$ws.Range("C3").Value = "10,500~14,000"
$ws.Range("D3").Value = "-"
Set-TextValue $ws.Range("E3") "13650"
$ws.Range("F3").Value = "한국투자증권"

# Row 4: 하스
$ws.Range("A4").Value = "하스"
$ws.Range("B4").Value = "2024.05.16~05.22"
$ws.Range("C4").Value = "9,000~12,000"
$ws.Range("D4").Value = "-"
Set-TextValue $ws.Range("E4") "16290"
$ws.Range("F4").Value = "삼성증권"

# Row 5: 미래에셋비전스팩4호
$ws.Range("A5").Value = "미래에셋비전스팩4호"
$ws.Range("B5").Value = "2024.05.13~05.14"
$ws.Range("C5").Value = "2,000~2,000"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = "미래에셋증권"

# Row 6: 노브랜드
$ws.Range("A6").Value = "노브랜드"
$ws.Range("B6").Value = "2024.04.30~05.08"
$ws.Range("C6").Value = "8,700~11,000"
$ws.Range("D6").Value = "-"
Set-TextValue $ws.Range("E6") "10440"
$ws.Range("F6").Value = "삼성증권"

# Row 7: 에스오에스랩
$ws.Range("A7").Value = "에스오에스랩"
$ws.Range("B7").Value = "2024.04.30~05.08"
$ws.Range("C7").Value = "7,500~9,000"
$ws.Range("D7").Value = "-"
Set-TextValue $ws.Range("E7") "15000"
$ws.Range("F7").Value = "한국투자증권"

# Row 8: KB스팩28호
$ws.Range("A8").Value = "KB스팩28호"
$ws.Range("B8").Value = "2024.04.29~04.30"
$ws.Range("C8").Value = "2,000~2,000"
$ws.Range("D8").Value = "-"
Set-TextValue $ws.Range("E8") "10000"
$ws.Range("F8").Value = "KB증권"

# Row 9: 아이씨티케이
$ws.Range("A9").Value = "아이씨티케이"
$ws.Range("B9").Value = "2024.04.24~04.30"
$ws.Range("C9").Value = "13,000~16,000"
$ws.Range("D9").Value = "-"
Set-TextValue $ws.Range("E9") "25610"
$ws.Range("F9").Value = "NH투자증권"

# Row 11: HD현대마린솔루션(구.HD현대글로벌서비스)(유가)
$ws.Range("A11").Value = "HD현대마린솔루션(구.HD현대글로벌서비스)(유가)"
$ws.Range("B11").Value = "2024.04.16~04.22"
$ws.Range("C11").Value = "73,300~83,400"
Set-TextValue $ws.Range("D11") "83400"
Set-TextValue $ws.Range("E11") "652370"
$ws.Range("F11").Value = "KB증권,신한투자증권,하나증권,대신증권,삼성증권"

Write-Output "RPA datasets push 2024-04-27 applied"